# The edit re-shuffles the observation rows 2-9 on the "Artfynd" sheet: each
# row's identifying/observation data (Id, Taxonsorteringsordning, Rödlistade,
# TaxonId, Artnamn, Vetenskapligt namn, Auktor, Lokalnamn, Ost, Nord) moves to
# a different row, per a fixed permutation; every other column already holds
# the same value in all eight rows, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ from row to row and therefore need to
# travel with the record when rows are reshuffled.
$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R")

# Snapshot the current (pre-edit) contents of those columns for rows 2-9
# before any writes happen, so later writes can't clobber a value we still
# need to read.
$snapshot = @{}
foreach ($col in $cols) {
    $rowVals = @{}
    for ($r = 2; $r -le 9; $r++) {
        $rowVals[$r] = $ws.Range("$col$r").Value2
    }
    $snapshot[$col] = $rowVals
}

# destinationRow -> sourceRow: the content that ends up in row N is whatever
# used to live in row Map[N].
$rowMap = @{
    2 = 4
    3 = 6
    4 = 3
    5 = 8
    6 = 2
    7 = 5
    8 = 9
    9 = 7
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot[$col][$srcRow]
    }
}
